$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: the two country labels in D1/E1 (HUN_GBR_841840 / NLD_GBR_841840)
# are reordered - swap which column they label.
$ws.Cells.Item(1, 4).Value2 = "NLD_GBR_841840"
$ws.Cells.Item(1, 5).Value2 = "HUN_GBR_841840"

# Data rows 2..171: columns D and E had their values swapped (re-ordering the
# two index series). Handle every row uniformly - this also naturally covers
# the rows where only one of D/E was populated (the empty side just moves
# across as an empty write, clearing the old cell and creating the new one).
for ($r = 2; $r -le 171; $r++) {
    $dVal = $ws.Cells.Item($r, 4).Value2
    $eVal = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 4).Value2 = $eVal
    $ws.Cells.Item($r, 5).Value2 = $dVal
}
